$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> MuSCs
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Cd80"
$ws.Cells.Item(2, 3).Value = "Ctla4"
$ws.Cells.Item(2, 4).Value = "MuSCs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 0.7715313333333333
$ws.Cells.Item(2, 8).Value = 2.314594
$ws.Cells.Item(2, 9).Value = 0.05172308417778351
$ws.Cells.Item(2, 10).Value = 0.05172308417778351
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.005070666666666667
$ws.Cells.Item(2, 14).Value = 0.015212
$ws.Cells.Item(2, 15).Value = 0.008865695667453655
$ws.Cells.Item(2, 16).Value = 0.008865695667453653
$ws.Cells.Item(2, 17).Value = 0.003912178214222222
$ws.Cells.Item(2, 18).Value = 0.035209603928
$ws.Cells.Item(2, 19).Value = 0.0004585611233023159
$ws.Cells.Item(2, 20).Value = 0.0004585611233023159

# Row 3: ECs -> Resolving-Mac
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Cd80"
$ws.Cells.Item(3, 3).Value = "Ctla4"
$ws.Cells.Item(3, 4).Value = "Resolving-Mac"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 0.7715313333333333
$ws.Cells.Item(3, 8).Value = 2.314594
$ws.Cells.Item(3, 9).Value = 0.05172308417778351
$ws.Cells.Item(3, 10).Value = 0.05172308417778351
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 0.5668716666666667
$ws.Cells.Item(3, 14).Value = 1.700615
$ws.Cells.Item(3, 15).Value = 0.9911343043325463
$ws.Cells.Item(3, 16).Value = 0.9911343043325463
$ws.Cells.Item(3, 17).Value = 0.4373592528122222
$ws.Cells.Item(3, 18).Value = 3.93623327531
$ws.Cells.Item(3, 19).Value = 0.05126452305448119
$ws.Cells.Item(3, 20).Value = 0.0512645230544812

# Row 4: FAPs -> MuSCs
$ws.Cells.Item(4, 1).Value = "FAPs"
$ws.Cells.Item(4, 2).Value = "Cd80"
$ws.Cells.Item(4, 3).Value = "Ctla4"
$ws.Cells.Item(4, 4).Value = "MuSCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 3.713472666666666
$ws.Cells.Item(4, 8).Value = 11.140418
$ws.Cells.Item(4, 9).Value = 0.2489493958723191
$ws.Cells.Item(4, 10).Value = 0.2489493958723191
$ws.Cells.Item(4, 11).Value = 1
$ws.Cells.Item(4, 12).Value = 0.3333333333333333
$ws.Cells.Item(4, 13).Value = 0.005070666666666667
$ws.Cells.Item(4, 14).Value = 0.015212
$ws.Cells.Item(4, 15).Value = 0.008865695667453655
$ws.Cells.Item(4, 16).Value = 0.008865695667453653
$ws.Cells.Item(4, 17).Value = 0.01882978206844444
$ws.Cells.Item(4, 18).Value = 0.169468038616
$ws.Cells.Item(4, 19).Value = 0.002207109580400424
$ws.Cells.Item(4, 20).Value = 0.002207109580400424

# Row 5: FAPs -> Resolving-Mac
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Cd80"
$ws.Cells.Item(5, 3).Value = "Ctla4"
$ws.Cells.Item(5, 4).Value = "Resolving-Mac"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 3.713472666666666
$ws.Cells.Item(5, 8).Value = 11.140418
$ws.Cells.Item(5, 9).Value = 0.2489493958723191
$ws.Cells.Item(5, 10).Value = 0.2489493958723191
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 0.5668716666666667
$ws.Cells.Item(5, 14).Value = 1.700615
$ws.Cells.Item(5, 15).Value = 0.9911343043325463
$ws.Cells.Item(5, 16).Value = 0.9911343043325463
$ws.Cells.Item(5, 17).Value = 2.105062439674444
$ws.Cells.Item(5, 18).Value = 18.94556195707
$ws.Cells.Item(5, 19).Value = 0.2467422862919186
$ws.Cells.Item(5, 20).Value = 0.2467422862919187

# Row 6: MuSCs -> MuSCs
$ws.Cells.Item(6, 1).Value = "MuSCs"
$ws.Cells.Item(6, 2).Value = "Cd80"
$ws.Cells.Item(6, 3).Value = "Ctla4"
$ws.Cells.Item(6, 4).Value = "MuSCs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 1.701035666666667
$ws.Cells.Item(6, 8).Value = 5.103107
$ws.Cells.Item(6, 9).Value = 0.1140366012048922
$ws.Cells.Item(6, 10).Value = 0.1140366012048922
$ws.Cells.Item(6, 11).Value = 1
$ws.Cells.Item(6, 12).Value = 0.3333333333333333
$ws.Cells.Item(6, 13).Value = 0.005070666666666667
$ws.Cells.Item(6, 14).Value = 0.015212
$ws.Cells.Item(6, 15).Value = 0.008865695667453655
$ws.Cells.Item(6, 16).Value = 0.008865695667453653
$ws.Cells.Item(6, 17).Value = 0.008625384853777777
$ws.Cells.Item(6, 18).Value = 0.07762846368399999
$ws.Cells.Item(6, 19).Value = 0.001011013801233353
$ws.Cells.Item(6, 20).Value = 0.001011013801233353

# Row 7: MuSCs -> Resolving-Mac
$ws.Cells.Item(7, 1).Value = "MuSCs"
$ws.Cells.Item(7, 2).Value = "Cd80"
$ws.Cells.Item(7, 3).Value = "Ctla4"
$ws.Cells.Item(7, 4).Value = "Resolving-Mac"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 1.701035666666667
$ws.Cells.Item(7, 8).Value = 5.103107
$ws.Cells.Item(7, 9).Value = 0.1140366012048922
$ws.Cells.Item(7, 10).Value = 0.1140366012048922
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 0.5668716666666667
$ws.Cells.Item(7, 14).Value = 1.700615
$ws.Cells.Item(7, 15).Value = 0.9911343043325463
$ws.Cells.Item(7, 16).Value = 0.9911343043325463
$ws.Cells.Item(7, 17).Value = 0.9642689234227777
$ws.Cells.Item(7, 18).Value = 8.678420310804999
$ws.Cells.Item(7, 19).Value = 0.1130255874036589
$ws.Cells.Item(7, 20).Value = 0.1130255874036589

# Row 8: Resolving-Mac -> MuSCs
$ws.Cells.Item(8, 1).Value = "Resolving-Mac"
$ws.Cells.Item(8, 2).Value = "Cd80"
$ws.Cells.Item(8, 3).Value = "Ctla4"
$ws.Cells.Item(8, 4).Value = "MuSCs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 8.730536666666668
$ws.Cells.Item(8, 8).Value = 26.19161
$ws.Cells.Item(8, 9).Value = 0.5852909187450052
$ws.Cells.Item(8, 10).Value = 0.5852909187450052
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 0.005070666666666667
$ws.Cells.Item(8, 14).Value = 0.015212
$ws.Cells.Item(8, 15).Value = 0.008865695667453655
$ws.Cells.Item(8, 16).Value = 0.008865695667453653
$ws.Cells.Item(8, 17).Value = 0.04426964125777778
$ws.Cells.Item(8, 18).Value = 0.39842677132
$ws.Cells.Item(8, 19).Value = 0.005189011162517561
$ws.Cells.Item(8, 20).Value = 0.00518901116251756

# Row 9: Resolving-Mac -> Resolving-Mac
$ws.Cells.Item(9, 1).Value = "Resolving-Mac"
$ws.Cells.Item(9, 2).Value = "Cd80"
$ws.Cells.Item(9, 3).Value = "Ctla4"
$ws.Cells.Item(9, 4).Value = "Resolving-Mac"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 8.730536666666668
$ws.Cells.Item(9, 8).Value = 26.19161
$ws.Cells.Item(9, 9).Value = 0.5852909187450052
$ws.Cells.Item(9, 10).Value = 0.5852909187450052
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 0.5668716666666667
$ws.Cells.Item(9, 14).Value = 1.700615
$ws.Cells.Item(9, 15).Value = 0.9911343043325463
$ws.Cells.Item(9, 16).Value = 0.9911343043325463
$ws.Cells.Item(9, 17).Value = 4.949093871127778
$ws.Cells.Item(9, 18).Value = 44.54184484015
$ws.Cells.Item(9, 19).Value = 0.5801019075824876
$ws.Cells.Item(9, 20).Value = 0.5801019075824876

